$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column C (y_0_forecast) and column E (y_1_forecast)
# Row number => (C value, E value)
$values = @{
    2  = @(1.392321641630434,   1.710071460977503)
    3  = @(1.004409005705997,   1.642433761320072)
    4  = @(2.418114148635109,   2.828066716168021)
    5  = @(1.89159218653383,    2.544631191216329)
    6  = @(1.554977796875501,   1.312870290004287)
    7  = @(0.6180254938795482,  0.7749619016293785)
    8  = @(0.481899667566732,   0.7487574275252262)
    9  = @(1.905862317202089,   1.389591155234515)
    10 = @(2.671046044496239,   2.125743999456575)
    11 = @(1.917627847674064,   2.694711744616662)
    12 = @(1.119562422009102,   1.831617848540201)
    13 = @(1.344920716048192,   1.037735724446631)
    14 = @(2.195375580740766,   1.872521508785896)
    15 = @(2.542856270410665,   2.961494745505977)
    16 = @(0.3979826440748235,  2.008592810942544)
    17 = @(-2.604000402888396, -0.08252516517808228)
    18 = @(1.122551915563408,   0.254631175783615)
    19 = @(2.273132718878146,   1.620205313802381)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}
